$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8 & 9, and add the new "raw working photos" rows 10-14 ---
# Order chosen to reproduce the shared-string append order seen in the target file.

# Row 9: Status changes from "In Queue" to "Scrubbed"
$ws.Range("D9").Value = "Scrubbed"

# Row 10 (new): Beach Rock, still "In Queue"
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Beach Rock"
$ws.Range("D10").Value = "In Queue"

# Row 8: Image column filled in
$ws.Range("C8").Value = "To many to count"

# Row 9: Image column filled in
$ws.Range("C9").Value = "N/A"

# Row 11 (new): Sail Boats
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Sail Boats"
$ws.Range("D11").Value = "In Queue"

# Row 12 (new): Sunset
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Sunset"
$ws.Range("D12").Value = "In Queue"

# Row 13 (new): Water bowl
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Water bowl"
$ws.Range("D13").Value = "In Queue"

# Row 14 (new): Tiny Waterfall
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Tiny Waterfall"
$ws.Range("D14").Value = "In Queue"

# Row 8: Status changes from "In Queue" to "Posted"
$ws.Range("D8").Value = "Posted"

# --- Data validation list on column D: extend range and add "Scrubbed" option ---
$ws.Range("D2:D30").Validation.Delete()
$ws.Range("D2:D30").Validation.Add(3)
$ws.Range("D2:D30").Validation.Formula1 = """In Queue, Edited, Written, Posted, Scrubbed"""

# --- Restore the active cell selection to match the saved view state ---
$ws.Range("F19").Select()
